# add guild data module
# Appends a new "GuilID" field row to the Property sheet of Guild.xlsx,
# describing the public/private/save/view flags, relation value and the
# (Chinese) description for the new field, mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "GuilID"
$ws.Range("B11").Value = "object"
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = $true
$ws.Range("F11").Value = $true
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = "Friend"
$ws.Range("J11").Value = "工会ID"

# Match the text-formatted style used by the other Id/Type/RelationValue/Desc cells
$ws.Range("A11").NumberFormat = "@"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("J11").NumberFormat = "@"

# Move the selection the way the author's workbook ended up (E19 instead of E20)
[void]$ws.Range("E19").Select()
